$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.050.30"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "3.065.16"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.38"
$ws.Range("E5").Value = "  -3.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "606.27"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.09"
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.376"
$ws.Range("E8").Value = "  -6.24%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  +5.77%  "
$ws.Range("D11").Value = "3.062.62"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.195"
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("D13").Value = "93.709.09"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000238"
$ws.Range("E14").Value = "  -5.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.47"
$ws.Range("E15").Value = "  -3.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.28"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").Value = "3.631.52"
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("D18").Value = "3.044.23"
$ws.Range("E18").Value = "  -2.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.52"
$ws.Range("E19").Value = "  -7.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.27"
$ws.Range("E20").Value = "  -4.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.65"
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "436.14"
$ws.Range("E22").Value = "  -2.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.73"
$ws.Range("E23").Value = "  -6.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000187"
$ws.Range("E24").Value = "  -9.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.33"
$ws.Range("E25").Value = "  +5.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.46"
$ws.Range("E26").Value = "  -6.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "84.03"
$ws.Range("E27").Value = "  -3.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.68"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").Value = "3.222.91"
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.241"
$ws.Range("E31").Value = "  +4.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.176"
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.122"
$ws.Range("E33").Value = "  -10.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").Value = "  +15.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.98"
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.59"
$ws.Range("E36").Value = "  -5.90%  "
$ws.Range("E37").Value = "  -3.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.21"
$ws.Range("E38").Value = "  -4.01%  "
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.98"
$ws.Range("E40").Value = "  +3.75%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.433"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("B42").Value = "MantraDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.71"
$ws.Range("E42").Value = "  -3.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "461.42"
$ws.Range("E43").Value = "  -6.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.24"
$ws.Range("E44").Value = "  -4.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.09"
$ws.Range("E46").Value = "  -9.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.70"
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.81"
$ws.Range("E48").Value = "  -5.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.662"
$ws.Range("E49").Value = "  -4.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.56"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.997"
$ws.Range("E51").Value = "  -0.23%  "
